# Insert a new weekly record at the top of the data table (row 33),
# shifting all existing records down by one row (old row 33 -> 34,
# old row 96 -> 97, etc. — same logical rows, just slotted one lower).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with this week's record.
$ws.Range("A33").Value = 5
$ws.Range("B33").Value = "Macroferia Regional de Talca"
$ws.Range("C33").Value = "Maule"
$ws.Range("D33").Value = 44915
$ws.Range("E33").Value = 7
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100101
$ws.Range("H33").Value = "Berries"
$ws.Range("I33").Value = 100101001
$ws.Range("J33").Value = "Arándano (blue)"
$ws.Range("K33").Value = "Sin especificar"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 220
$ws.Range("N33").Value = 3000
$ws.Range("O33").Value = 3200
$ws.Range("P33").Value = 3109
$ws.Range("Q33").Value = "$/bandeja 2 kilos"
$ws.Range("R33").Value = "Provincia de Curicó"
$ws.Range("S33").Value = 1554
$ws.Range("T33").Value = 2
